# Auto-generated edit script applying scheduled-runner price updates
# to Seraph_Profits.xlsx (market-data refresh across ALC/ARM/BSM/CRP/CUL/GSM/LTW/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 34.705883
$ws.Range("I11").Value = 34.705883
$ws.Range("K11").Value = 34.705883
$ws.Range("M11").Value = 105.294117

# Row 31
$ws.Range("H31").Value = 108.44444
$ws.Range("I31").Value = 108.44444
$ws.Range("K31").Value = 325.33332
$ws.Range("M31").Value = -95.33332000000001

# Row 40
$ws.Range("H40").Value = 2555.2222
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1825

# Row 55
$ws.Range("H55").Value = 227.04762
$ws.Range("I55").Value = 223.45
$ws.Range("J55").Value = 299
$ws.Range("K55").Value = 223.45
$ws.Range("L55").Value = 299
$ws.Range("M55").Value = -9.449999999999989
$ws.Range("N55").Value = -727

# Row 64
$ws.Range("H64").Value = 3549.5
$ws.Range("I64").Value = 3399.3333
$ws.Range("K64").Value = 3399.3333
$ws.Range("M64").Value = -3151.3333

# Row 67
$ws.Range("H67").Value = 3549.5
$ws.Range("I67").Value = 3399.3333
$ws.Range("K67").Value = 3399.3333
$ws.Range("M67").Value = -2541.3333

# Row 86
$ws.Range("H86").Value = 550002
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 550002
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 550002
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -552248

# Row 89
$ws.Range("H89").Value = 550002
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 550002
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 2750010
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -2761242

# Row 98
$ws.Range("H98").Value = 1690.75
$ws.Range("I98").Value = 1709.1
$ws.Range("K98").Value = 1709.1
$ws.Range("M98").Value = -211.0999999999999

# Row 100
$ws.Range("H100").Value = 1072.7273
$ws.Range("I100").Value = 1080
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 1080
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -539
$ws.Range("N100").Value = -2082

# Row 113
$ws.Range("H113").Value = 8810.888999999999
$ws.Range("I113").Value = 8659.799999999999
$ws.Range("K113").Value = 8659.799999999999
$ws.Range("M113").Value = -5405.799999999999

# Row 122
$ws.Range("H122").Value = 1690.75
$ws.Range("I122").Value = 1709.1
$ws.Range("K122").Value = 5127.299999999999
$ws.Range("M122").Value = -2677.299999999999

# Row 137
$ws.Range("H137").Value = 1855
$ws.Range("I137").Value = 1137.1
$ws.Range("K137").Value = 3411.3
$ws.Range("M137").Value = -861.2999999999997

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 7537.8184
$ws.Range("I61").Value = 7791.6
$ws.Range("K61").Value = 7791.6
$ws.Range("M61").Value = -7579.6

# Row 97
$ws.Range("H97").Value = 542.3333
$ws.Range("I97").Value = 406.83334
$ws.Range("K97").Value = 406.83334
$ws.Range("M97").Value = 89.16665999999998

# Row 136
$ws.Range("H136").Value = 7537.8184
$ws.Range("I136").Value = 7791.6
$ws.Range("K136").Value = 23374.8
$ws.Range("M136").Value = -20824.8

$ws = $wb.Worksheets.Item("BSM")
# Row 33
$ws.Range("H33").Value = 15500
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

# Row 76
$ws.Range("H76").Value = 29999
$ws.Range("J76").Value = 29999
$ws.Range("L76").Value = 29999
$ws.Range("N76").Value = -30629

# Row 79
$ws.Range("H79").Value = 29999
$ws.Range("J79").Value = 29999
$ws.Range("L79").Value = 29999
$ws.Range("N79").Value = -32183

# Row 99
$ws.Range("H99").Value = 1936.6666
$ws.Range("I99").Value = 2071.8333
$ws.Range("K99").Value = 2071.8333
$ws.Range("M99").Value = -573.8332999999998

# Row 134
$ws.Range("H134").Value = 3262.75
$ws.Range("I134").Value = 2565.3
$ws.Range("K134").Value = 7695.900000000001
$ws.Range("M134").Value = -5160.900000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 55558692
$ws.Range("I16").Value = 55558692
$ws.Range("K16").Value = 55558692
$ws.Range("M16").Value = -55558405

# Row 22
$ws.Range("H22").Value = 67565.95
$ws.Range("I22").Value = 86272.71000000001
$ws.Range("K22").Value = 86272.71000000001
$ws.Range("M22").Value = -85922.71000000001

# Row 31
$ws.Range("H31").Value = 2562.1035
$ws.Range("I31").Value = 1579.3914
$ws.Range("K31").Value = 1579.3914
$ws.Range("M31").Value = -1284.3914

# Row 34
$ws.Range("H34").Value = 2562.1035
$ws.Range("I34").Value = 1579.3914
$ws.Range("K34").Value = 1579.3914
$ws.Range("M34").Value = -1377.3914

# Row 54
$ws.Range("H54").Value = 18000
$ws.Range("J54").Value = 18000
$ws.Range("L54").Value = 18000
$ws.Range("N54").Value = -19316

# Row 99
$ws.Range("H99").Value = 2627.5
$ws.Range("I99").Value = 2627.5
$ws.Range("K99").Value = 2627.5
$ws.Range("M99").Value = -1129.5

# Row 110
$ws.Range("H110").Value = 80000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180

# Row 113
$ws.Range("H113").Value = 55558692
$ws.Range("I113").Value = 55558692
$ws.Range("K113").Value = 55558692
$ws.Range("M113").Value = -55556522

# Row 126
$ws.Range("H126").Value = 2627.5
$ws.Range("I126").Value = 2627.5
$ws.Range("K126").Value = 7882.5
$ws.Range("M126").Value = -5412.5

# Row 132
$ws.Range("H132").Value = 2012.4117
$ws.Range("I132").Value = 1872.2858
$ws.Range("K132").Value = 5616.857400000001
$ws.Range("M132").Value = -3086.857400000001

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 403.9
$ws.Range("I2").Value = 67.5
$ws.Range("J2").Value = 488
$ws.Range("K2").Value = 405
$ws.Range("L2").Value = 2928
$ws.Range("M2").Value = -292
$ws.Range("N2").Value = -3154

# Row 40
$ws.Range("H40").Value = 48.272728
$ws.Range("I40").Value = 29.25
$ws.Range("K40").Value = 117
$ws.Range("M40").Value = -48

# Row 131
$ws.Range("H131").Value = 1368.0444
$ws.Range("J131").Value = 1388.05
$ws.Range("L131").Value = 4164.15
$ws.Range("N131").Value = -14244.15

# Row 137
$ws.Range("H137").Value = 5999.8
$ws.Range("I137").Value = 6666.3335
$ws.Range("K137").Value = 19999.0005
$ws.Range("M137").Value = -14899.0005

# Row 139
$ws.Range("H139").Value = 1535.3334
$ws.Range("I139").Value = 1535.3334
$ws.Range("K139").Value = 4606.0002
$ws.Range("M139").Value = 533.9997999999996

# Row 140
$ws.Range("H140").Value = 1068.4
$ws.Range("I140").Value = 861.5789
$ws.Range("K140").Value = 2584.7367
$ws.Range("M140").Value = 2595.2633

$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Range("H41").Value = 7499.5
$ws.Range("I41").Value = 7499.5
$ws.Range("K41").Value = 7499.5
$ws.Range("M41").Value = -7144.5

# Row 80
$ws.Range("H80").Value = 5247.25
$ws.Range("I80").Value = 4996.3335
$ws.Range("J80").Value = 6000
$ws.Range("K80").Value = 4996.3335
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -3998.3335
$ws.Range("N80").Value = -7996

# Row 83
$ws.Range("H83").Value = 5247.25
$ws.Range("I83").Value = 4996.3335
$ws.Range("J83").Value = 6000
$ws.Range("K83").Value = 24981.6675
$ws.Range("L83").Value = 30000
$ws.Range("M83").Value = -19989.6675
$ws.Range("N83").Value = -39984

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 6301.6553
$ws.Range("I22").Value = 2429.8667
$ws.Range("K22").Value = 2429.8667
$ws.Range("M22").Value = -2134.8667

# Row 27
$ws.Range("H27").Value = 6301.6553
$ws.Range("I27").Value = 2429.8667
$ws.Range("K27").Value = 2429.8667
$ws.Range("M27").Value = -2322.8667

# Row 40
$ws.Range("H40").Value = 3864.8333
$ws.Range("I40").Value = 3864.8333
$ws.Range("K40").Value = 3864.8333
$ws.Range("M40").Value = -3728.8333

# Row 61
$ws.Range("H61").Value = 13892064
$ws.Range("I61").Value = 15876216
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 15876216
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -15876014
$ws.Range("N61").Value = -3404

# Row 82
$ws.Range("H82").Value = 144690.14
$ws.Range("I82").Value = 959.25
$ws.Range("K82").Value = 959.25
$ws.Range("M82").Value = -598.25

# Row 85
$ws.Range("H85").Value = 144690.14
$ws.Range("I85").Value = 959.25
$ws.Range("K85").Value = 959.25
$ws.Range("M85").Value = 288.75

# Row 100
$ws.Range("H100").Value = 6948
$ws.Range("I100").Value = 6877.3
$ws.Range("J100").Value = 7301.5
$ws.Range("K100").Value = 6877.3
$ws.Range("L100").Value = 7301.5
$ws.Range("M100").Value = -6336.3
$ws.Range("N100").Value = -8383.5

# Row 113
$ws.Range("H113").Value = 13892064
$ws.Range("I113").Value = 15876216
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 15876216
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -15874046
$ws.Range("N113").Value = -7340

# Row 122
$ws.Range("H122").Value = 6841.1763
$ws.Range("I122").Value = 6333.75
$ws.Range("J122").Value = 6997.3076
$ws.Range("K122").Value = 19001.25
$ws.Range("L122").Value = 20991.9228
$ws.Range("M122").Value = -16551.25
$ws.Range("N122").Value = -25891.9228

# Row 132
$ws.Range("H132").Value = 73900.36
$ws.Range("I132").Value = 85633.75
$ws.Range("K132").Value = 256901.25
$ws.Range("M132").Value = -254371.25

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2959.12
$ws.Range("I122").Value = 2728.4
$ws.Range("J122").Value = 3305.2
$ws.Range("K122").Value = 8185.200000000001
$ws.Range("L122").Value = 9915.599999999999
$ws.Range("M122").Value = -5735.200000000001
$ws.Range("N122").Value = -14815.6
